# Generate Report for Handback
#
# Row 6 (e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e) has now been handed back and is
# in sync with en-US; row 7 (0b52ffe7-7526-47a6-a9b1-f913f9557407) is still
# pending (transform failed). Swap their positions on every sheet (Overview,
# zh-cn, de-de) and fill in e0f18d0f's now-completed handback columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A6").Value = "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.md"
$ov.Range("B6").Value = "Handed back: in sync with en-US"
$ov.Range("C6").Value = "Handed back: in sync with en-US"
$ov.Range("D6").Value = "2016-03-24 15:15:31"

$ov.Range("A7").Value = "0b52ffe7-7526-47a6-a9b1-f913f9557407.md"
$ov.Range("B7").Value = "Handback transform failed"
$ov.Range("C7").Value = "Handback transform failed"
$ov.Range("D7").Value = "2016-03-24 15:07:41"

foreach ($h in $ov.Hyperlinks) {
    $rref = $h.Range.Address()
    if ($rref -eq '$A$6') {
        $h.Address = "https://github.com/OpenLocalizationTest/oltest/blob/ec56cec02dbdd6527f76f8aca445624cdf219bca/e2e/e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.md"
        $h.TextToDisplay = "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.md"
    }
    elseif ($rref -eq '$A$7') {
        $h.Address = "https://github.com/OpenLocalizationTest/oltest/blob/94f848faba723f7059e59111dfdc1d0ccccd47bc/e2e/0b52ffe7-7526-47a6-a9b1-f913f9557407.md"
        $h.TextToDisplay = "0b52ffe7-7526-47a6-a9b1-f913f9557407.md"
    }
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 6 becomes e0f18d0f, now fully handed back (Latest Target File / Latest
# Handback File columns get populated to mirror the handoff columns).
$zh.Range("A6").Value = "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.md"
$zh.Range("C6").Value = "Handed back: in sync with en-US"
$zh.Range("D6").Value = "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.c241a454f2fa57108d2f88b9821e8577c52b0657.zh-cn.xlf"
$zh.Range("E6").Value = "2016-03-24 15:15:26"
$zh.Range("F6").Value = "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.md"
$zh.Range("G6").Value = "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.c241a454f2fa57108d2f88b9821e8577c52b0657.zh-cn.xlf"
$zh.Range("H6").Value = "2016-03-24 15:17:06"
$zh.Range("L6").Value = ""

# Row 7 becomes 0b52ffe7, still failing, same content as the old row 6 (minus
# the columns that were never populated for it).
$zh.Range("A7").Value = "0b52ffe7-7526-47a6-a9b1-f913f9557407.md"
$zh.Range("C7").Value = "Handback transform failed"
$zh.Range("D7").Value = "0b52ffe7-7526-47a6-a9b1-f913f9557407.925700ff9c9abff613f6ca1542b91f2c745de3dc.zh-cn.xlf"
$zh.Range("E7").Value = "2016-03-24 15:07:36"
$zh.Range("F7").Value = ""
$zh.Range("G7").Value = ""
$zh.Range("H7").Value = "0001-01-01 00:00:00"
$zh.Range("L7").Value = "Handback file name: x2kyg3ye.yxh is different with handoff file name: 0b52ffe7-7526-47a6-a9b1-f913f9557407.925700ff9c9abff613f6ca1542b91f2c745de3dc.zh-cn."

foreach ($h in $zh.Hyperlinks) {
    $rref = $h.Range.Address()
    if ($rref -eq '$A$6') {
        $h.Address = "https://github.com/OpenLocalizationTest/oltest/blob/ec56cec02dbdd6527f76f8aca445624cdf219bca/e2e/e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.md"
        $h.TextToDisplay = "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.md"
    }
    elseif ($rref -eq '$D$6') {
        $h.Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/07e397b9c626ca08b4755d1f039a6f37e7c26129/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.c241a454f2fa57108d2f88b9821e8577c52b0657.zh-cn.xlf"
        $h.TextToDisplay = "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.c241a454f2fa57108d2f88b9821e8577c52b0657.zh-cn.xlf"
    }
    elseif ($rref -eq '$A$7') {
        $h.Address = "https://github.com/OpenLocalizationTest/oltest/blob/94f848faba723f7059e59111dfdc1d0ccccd47bc/e2e/0b52ffe7-7526-47a6-a9b1-f913f9557407.md"
        $h.TextToDisplay = "0b52ffe7-7526-47a6-a9b1-f913f9557407.md"
    }
    elseif ($rref -eq '$D$7') {
        $h.Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/965b05319fcacd79ccd11e8b3a13d2e124d35972/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0b52ffe7-7526-47a6-a9b1-f913f9557407.925700ff9c9abff613f6ca1542b91f2c745de3dc.zh-cn.xlf"
        $h.TextToDisplay = "0b52ffe7-7526-47a6-a9b1-f913f9557407.925700ff9c9abff613f6ca1542b91f2c745de3dc.zh-cn.xlf"
    }
}

$zh.Hyperlinks.Add($zh.Range("F6"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e0f18d0f1c534dc58c3ee5429e9ec73e0000001/e2e/e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.md", "", "", "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("G6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e0f18d0f1c534dc58c3ee5429e9ec73e0000002/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.c241a454f2fa57108d2f88b9821e8577c52b0657.zh-cn.xlf", "", "", "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.c241a454f2fa57108d2f88b9821e8577c52b0657.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A6").Value = "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.md"
$de.Range("C6").Value = "Handed back: in sync with en-US"
$de.Range("E6").Value = "2016-03-24 15:15:31"
$de.Range("F6").Value = "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.md"
$de.Range("G6").Value = "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.c241a454f2fa57108d2f88b9821e8577c52b0657.de-de.xlf"
$de.Range("H6").Value = "2016-03-24 15:17:13"
$de.Range("L6").Value = ""

$de.Range("A7").Value = "0b52ffe7-7526-47a6-a9b1-f913f9557407.md"
$de.Range("C7").Value = "Handback transform failed"
$de.Range("E7").Value = "2016-03-24 15:07:41"
$de.Range("F7").Value = ""
$de.Range("G7").Value = ""
$de.Range("H7").Value = "0001-01-01 00:00:00"
$de.Range("L7").Value = "Handback file name: x2kyg3ye.yxh is different with handoff file name: 0b52ffe7-7526-47a6-a9b1-f913f9557407.925700ff9c9abff613f6ca1542b91f2c745de3dc.de-de."

foreach ($h in $de.Hyperlinks) {
    $rref = $h.Range.Address()
    if ($rref -eq '$A$6') {
        $h.Address = "https://github.com/OpenLocalizationTest/oltest/blob/ec56cec02dbdd6527f76f8aca445624cdf219bca/e2e/e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.md"
        $h.TextToDisplay = "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.md"
    }
    elseif ($rref -eq '$D$6') {
        $h.Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/07e397b9c626ca08b4755d1f039a6f37e7c26129/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.c241a454f2fa57108d2f88b9821e8577c52b0657.de-de.xlf"
        $h.TextToDisplay = "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.c241a454f2fa57108d2f88b9821e8577c52b0657.de-de.xlf"
    }
    elseif ($rref -eq '$A$7') {
        $h.Address = "https://github.com/OpenLocalizationTest/oltest/blob/94f848faba723f7059e59111dfdc1d0ccccd47bc/e2e/0b52ffe7-7526-47a6-a9b1-f913f9557407.md"
        $h.TextToDisplay = "0b52ffe7-7526-47a6-a9b1-f913f9557407.md"
    }
    elseif ($rref -eq '$D$7') {
        $h.Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/965b05319fcacd79ccd11e8b3a13d2e124d35972/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0b52ffe7-7526-47a6-a9b1-f913f9557407.925700ff9c9abff613f6ca1542b91f2c745de3dc.de-de.xlf"
        $h.TextToDisplay = "0b52ffe7-7526-47a6-a9b1-f913f9557407.925700ff9c9abff613f6ca1542b91f2c745de3dc.de-de.xlf"
    }
}

$de.Hyperlinks.Add($de.Range("F6"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/e0f18d0f1c534dc58c3ee5429e9ec73e0000003/e2e/e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.md", "", "", "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.md") | Out-Null
$de.Hyperlinks.Add($de.Range("G6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e0f18d0f1c534dc58c3ee5429e9ec73e0000004/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.c241a454f2fa57108d2f88b9821e8577c52b0657.de-de.xlf", "", "", "e0f18d0f-1c53-4dc5-8c3e-e5429e9ec73e.c241a454f2fa57108d2f88b9821e8577c52b0657.de-de.xlf") | Out-Null
